$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 952.44775
$ws.Range("I129").Value = 520.8182
$ws.Range("J129").Value = 1037.2322
$ws.Range("K129").Value = 1562.4546
$ws.Range("L129").Value = 3111.6966
$ws.Range("M129").Value = 3437.5454
$ws.Range("N129").Value = -13111.6966
$ws.Range("H138").Value = 2529.7297
$ws.Range("I138").Value = 1663.8518
$ws.Range("J138").Value = 4867.6
$ws.Range("K138").Value = 4991.555399999999
$ws.Range("L138").Value = 14602.8
$ws.Range("M138").Value = 148.4446000000007
$ws.Range("N138").Value = -24882.8
$ws.Range("H141").Value = 5834.1875
$ws.Range("I141").Value = 2741
$ws.Range("J141").Value = 35735
$ws.Range("K141").Value = 8223
$ws.Range("L141").Value = 107205
$ws.Range("M141").Value = -3043
$ws.Range("N141").Value = -117565

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 20000
$ws.Range("I37").Value = 20000
$ws.Range("K37").Value = 20000
$ws.Range("M37").Value = -19727
$ws.Range("H74").Value = 996.5172
$ws.Range("I74").Value = 934.88
$ws.Range("J74").Value = 1381.75
$ws.Range("K74").Value = 934.88
$ws.Range("L74").Value = 1381.75
$ws.Range("M74").Value = -60.88
$ws.Range("N74").Value = -3129.75
$ws.Range("H77").Value = 996.5172
$ws.Range("I77").Value = 934.88
$ws.Range("J77").Value = 1381.75
$ws.Range("K77").Value = 4674.4
$ws.Range("L77").Value = 6908.75
$ws.Range("M77").Value = -306.3999999999996
$ws.Range("N77").Value = -15644.75
$ws.Range("H123").Value = 27309.111
$ws.Range("J123").Value = 27309.111
$ws.Range("L123").Value = 27309.111
$ws.Range("N123").Value = -37109.111
$ws.Range("H132").Value = 1889.6364
$ws.Range("I132").Value = 1190.2693
$ws.Range("J132").Value = 2899.8333
$ws.Range("K132").Value = 3570.8079
$ws.Range("L132").Value = 8699.499899999999
$ws.Range("M132").Value = -1040.8079
$ws.Range("N132").Value = -13759.4999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 40074
$ws.Range("J35").Value = 40074
$ws.Range("L35").Value = 40074
$ws.Range("N35").Value = -40694

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 2087.3784
$ws.Range("I31").Value = 1475
$ws.Range("K31").Value = 1475
$ws.Range("M31").Value = -1180
$ws.Range("H34").Value = 2087.3784
$ws.Range("I34").Value = 1475
$ws.Range("K34").Value = 1475
$ws.Range("M34").Value = -1273
$ws.Range("H58").Value = 1196341.2
$ws.Range("I58").Value = 2470734
$ws.Range("J58").Value = 1598.125
$ws.Range("K58").Value = 2470734
$ws.Range("L58").Value = 1598.125
$ws.Range("M58").Value = -2470531
$ws.Range("N58").Value = -2004.125
$ws.Range("H99").Value = 7125
$ws.Range("I99").Value = 7125
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 7125
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -5627
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 8417.786
$ws.Range("I105").Value = 8988.385
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 8988.385
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = -7241.385
$ws.Range("N105").Value = -4494
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 7125
$ws.Range("I126").Value = 7125
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 21375
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -18905
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 411460.66
$ws.Range("I132").Value = 564362.3
$ws.Range("J132").Value = 3723
$ws.Range("K132").Value = 1693086.9
$ws.Range("L132").Value = 11169
$ws.Range("M132").Value = -1690556.9
$ws.Range("N132").Value = -16229
$ws.Range("H134").Value = 2256.6785
$ws.Range("I134").Value = 1770.7646
$ws.Range("J134").Value = 3007.6365
$ws.Range("K134").Value = 5312.293799999999
$ws.Range("L134").Value = 9022.9095
$ws.Range("M134").Value = -2777.293799999999
$ws.Range("N134").Value = -14092.9095
$ws.Range("H136").Value = 1196341.2
$ws.Range("I136").Value = 2470734
$ws.Range("J136").Value = 1598.125
$ws.Range("K136").Value = 7412202
$ws.Range("L136").Value = 4794.375
$ws.Range("M136").Value = -7409652
$ws.Range("N136").Value = -9894.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("H131").Value = 20856034
$ws.Range("I131").Value = 12470
$ws.Range("J131").Value = 25666088
$ws.Range("K131").Value = 37410
$ws.Range("L131").Value = 76998264
$ws.Range("M131").Value = -32370
$ws.Range("N131").Value = -77008344

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 34204040
$ws.Range("J43").Value = 5733.3335
$ws.Range("L43").Value = 5733.3335
$ws.Range("N43").Value = -6035.3335
$ws.Range("H46").Value = 16425
$ws.Range("J46").Value = 19900
$ws.Range("L46").Value = 19900
$ws.Range("N46").Value = -20212
$ws.Range("H97").Value = 108677.14
$ws.Range("I97").Value = 47270.91
$ws.Range("J97").Value = 333833.34
$ws.Range("K97").Value = 47270.91
$ws.Range("L97").Value = 333833.34
$ws.Range("M97").Value = -46774.91
$ws.Range("N97").Value = -334825.34
$ws.Range("H116").Value = 45000
$ws.Range("J116").Value = 45000
$ws.Range("L116").Value = 45000
$ws.Range("N116").Value = -54178
$ws.Range("H122").Value = 3535.375
$ws.Range("I122").Value = 3099.25
$ws.Range("J122").Value = 3971.5
$ws.Range("K122").Value = 9297.75
$ws.Range("L122").Value = 11914.5
$ws.Range("M122").Value = -6847.75
$ws.Range("N122").Value = -16814.5
$ws.Range("H123").Value = 30322
$ws.Range("J123").Value = 30322
$ws.Range("L123").Value = 30322
$ws.Range("N123").Value = -35222
$ws.Range("H132").Value = 1618.45
$ws.Range("I132").Value = 960.875
$ws.Range("J132").Value = 4248.75
$ws.Range("K132").Value = 2882.625
$ws.Range("L132").Value = 12746.25
$ws.Range("M132").Value = -352.625
$ws.Range("N132").Value = -17806.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6177.6665
$ws.Range("I7").Value = 6449.75
$ws.Range("J7").Value = 5960
$ws.Range("K7").Value = 6449.75
$ws.Range("L7").Value = 5960
$ws.Range("M7").Value = -6337.75
$ws.Range("N7").Value = -6184
$ws.Range("H22").Value = 1180.1818
$ws.Range("J22").Value = 1422.75
$ws.Range("L22").Value = 1422.75
$ws.Range("N22").Value = -2012.75
$ws.Range("H27").Value = 1180.1818
$ws.Range("J27").Value = 1422.75
$ws.Range("L27").Value = 1422.75
$ws.Range("N27").Value = -1636.75
$ws.Range("H40").Value = 2178.7058
$ws.Range("I40").Value = 1772.1538
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 1772.1538
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -1636.1538
$ws.Range("N40").Value = -3772
$ws.Range("H126").Value = 6177.6665
$ws.Range("I126").Value = 6449.75
$ws.Range("J126").Value = 5960
$ws.Range("K126").Value = 19349.25
$ws.Range("L126").Value = 17880
$ws.Range("M126").Value = -16879.25
$ws.Range("N126").Value = -22820
$ws.Range("H132").Value = 5717.091
$ws.Range("I132").Value = 5292.4614
$ws.Range("J132").Value = 6330.4443
$ws.Range("K132").Value = 15877.3842
$ws.Range("L132").Value = 18991.3329
$ws.Range("M132").Value = -13347.3842
$ws.Range("N132").Value = -24051.3329
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 2020.2963
$ws.Range("I136").Value = 1805.1025
$ws.Range("J136").Value = 2579.8
$ws.Range("K136").Value = 5415.3075
$ws.Range("L136").Value = 7739.400000000001
$ws.Range("M136").Value = -2865.3075
$ws.Range("N136").Value = -12839.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 6893.5
$ws.Range("I100").Value = 14036.934
$ws.Range("J100").Value = 1253.9474
$ws.Range("K100").Value = 28073.868
$ws.Range("L100").Value = 2507.8948
$ws.Range("M100").Value = -27532.868
$ws.Range("N100").Value = -3589.8948
$ws.Range("H122").Value = 14882366
$ws.Range("I122").Value = 13890171
$ws.Range("J122").Value = 20835536
$ws.Range("K122").Value = 41670513
$ws.Range("L122").Value = 62506608
$ws.Range("M122").Value = -41668063
$ws.Range("N122").Value = -62511508
$ws.Range("H126").Value = 7688.75
$ws.Range("I126").Value = 8827.5
$ws.Range("J126").Value = 1995
$ws.Range("K126").Value = 26482.5
$ws.Range("L126").Value = 5985
$ws.Range("M126").Value = -24012.5
$ws.Range("N126").Value = -10925
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 1335.725
$ws.Range("I132").Value = 910.0345
$ws.Range("J132").Value = 2458
$ws.Range("K132").Value = 2730.1035
$ws.Range("L132").Value = 7374
$ws.Range("M132").Value = -200.1035000000002
$ws.Range("N132").Value = -12434

